# Updates cryptocurrency price/volume data in the worksheet to reflect
# the latest scrape, as produced by the GitHub Actions workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.556.77"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").Value = "2.643.98"
$ws.Range("E3").Value = "  -3.53%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'599.08"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "'167.88"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.544"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "2.642.41"
$ws.Range("E9").Value = "  -3.55%  "
$ws.Range("D10").Value = "'0.146"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").Value = "'0.159"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").Value = "'0.365"
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("D13").Value = "'5.24"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("D14").Value = "'27.99"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("D15").Value = "3.118.84"
$ws.Range("E15").Value = "  -3.68%  "
$ws.Range("D16").Value = "'0.0000185"
$ws.Range("E16").Value = "  -4.16%  "
$ws.Range("D17").Value = "67.407.98"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").Value = "2.636.67"
$ws.Range("E18").Value = "  -2.97%  "
$ws.Range("D19").Value = "'11.90"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").Value = "'7.87"
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("D21").Value = "'363.41"
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("D22").Value = "'4.41"
$ws.Range("E22").Value = "  -3.34%  "
$ws.Range("D23").Value = "'4.79"
$ws.Range("E23").Value = "  -4.55%  "
$ws.Range("D24").Value = "'10.85"
$ws.Range("E24").Value = "  +7.48%  "
$ws.Range("E25").Value = "  -5.47%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'70.76"
$ws.Range("E27").Value = "  -4.33%  "
$ws.Range("D28").Value = "2.779.38"
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("D29").Value = "'0.0000103"
$ws.Range("E29").Value = "  -4.52%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "'556.59"
$ws.Range("E31").Value = "  -5.45%  "
$ws.Range("D32").Value = "'8.05"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("D33").Value = "'1.40"
$ws.Range("E33").Value = "  -4.58%  "
$ws.Range("D34").Value = "'1.93"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -5.98%  "
$ws.Range("D38").Value = "'158.03"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "'19.39"
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("D40").Value = "'0.373"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("E41").Value = "  -5.20%  "
$ws.Range("D42").Value = "'5.28"
$ws.Range("E42").Value = "  -4.58%  "
$ws.Range("D43").Value = "'17.94"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'2.53"
$ws.Range("E44").Value = "  -6.53%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'40.18"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("E47").Value = "  -4.45%  "
$ws.Range("D48").Value = "'0.596"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "'153.62"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").Value = "'3.89"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("D51").Value = "'1.73"
$ws.Range("E51").Value = "  -4.25%  "
